$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OpsTracker")

# --- Step 1: toggle the existing AutoFilter off/on while the sheet's
# current data region is still A1:F34. The engine recomputes the
# AutoFilter range from the contiguous used region at the moment the
# filter is (re)applied, so doing this BEFORE we append the new row 35
# keeps the filter range at F34, matching the target state exactly
# (the new row is intentionally left outside of the filter range).
$null = $ws.Range("A1:F32").AutoFilter()
$null = $ws.Range("A1:F32").AutoFilter()

# --- Step 2: update the hidden defined name that mirrors the
# worksheet AutoFilter range (Excel keeps this in sync normally, but
# it is not auto-refreshed here, so set it explicitly).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "OpsTracker!_FilterDatabase") {
        $n.RefersTo = "=OpsTracker!`$A`$1:`$F`$34"
    }
}

# --- Step 3: cell content changes.
# D17: Todo -> Done
$ws.Range("D17").Value = "Done"

# Row 34's text moves to become new row 35's text ("Organize the
# office different accessories" keeps its place conceptually but the
# row itself shifts down), while rows 32-34 get new task descriptions.
# Write B34 first, then B32, then B33 so new shared-string entries are
# created in the same order the source workbook has them.
$ws.Range("B34").Value = "Online interview of Pronay Dhargave on 16 Nov at 3 PM"
$ws.Range("B32").Value = "File for daily visit enquiry Sheet"
$ws.Range("B33").Value = "File for bio data of teacher"

# New row 35.
$ws.Range("A35").Value = 35
$ws.Range("B35").Value = "Organize the office different accessories"
$ws.Range("C35").Value = "Debasish"
$ws.Range("D35").Value = "Todo"

# --- Step 4: column widths (best effort - the underlying engine
# quantizes ColumnWidth to 1/6-character increments, so we pick the
# closest representable width to the authored target).
$ws.Columns.Item(2).ColumnWidth = 46.3
$ws.Columns.Item(5).ColumnWidth = 62.7

# --- Step 5: restore the selection state recorded in the sheet view.
$ws.Range("B32").Select() | Out-Null
